# VLAB #45838 fix bug in calc of settlement age
# Add new Stock Synthesis change-log entries for versions 3.30.10.00 / 3.30.10.01
# and update the "latest revision" banner cell (F1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 44: 3.30.10.00 / new / fleets / provide controls for bycatch fleets ---
$ws.Range("F44").Value = "provide controls for bycatch fleets"

# --- Row 45: 3.30.10.00 / new / forecast / F0.1 option ---
$ws.Range("F45").Value = "provide benchmark and forecast option to use F0.1; this is either/or with F(Btgt)"

# --- Row 46: 3.30.10.00 / new / forecast / SSB_virgin vs SSB_unfished clarification ---
$ws.Range("F46").Value = "clarify internal usage and output for SSB_virgin vs SSB_unfished (benchmark); add description to manual"

# --- Row 47: 3.30.10.00 / new / forecast / clarify forecast-report.sso output ---
$ws.Range("F47").Value = "clarify output in forecast-report.sso and SPR/YPR profile"

# Version numbers for rows 44-47 (all 3.30.10.00) -- set after the descriptions above
# so "3.30.10.00" is appended to the shared-string table after the four new
# description strings, matching the authoring order.
$ws.Range("B44").Value = "3.30.10.00"
$ws.Range("B45").Value = "3.30.10.00"
$ws.Range("B46").Value = "3.30.10.00"
$ws.Range("B47").Value = "3.30.10.00"

# Dates for rows 44-47
$ws.Range("A44").Value = 43066
$ws.Range("A45").Value = 43066
$ws.Range("A46").Value = 43066
$ws.Range("A47").Value = 43109

# Type / Category for rows 44-47
$ws.Range("C44").Value = "new"
$ws.Range("D44").Value = "fleets"
$ws.Range("G44").Value = "Yes"

$ws.Range("C45").Value = "new"
$ws.Range("D45").Value = "forecast"
$ws.Range("G45").Value = "No"

$ws.Range("C46").Value = "new"
$ws.Range("D46").Value = "forecast"
$ws.Range("G46").Value = "No"

$ws.Range("C47").Value = "new"
$ws.Range("D47").Value = "forecast"
$ws.Range("G47").Value = "No"

# Row 46 wraps onto two lines in the real workbook -- taller row height.
$ws.Rows.Item(46).RowHeight = 31.5

# Update the banner cell (F1) to reference the newest version/date; this also
# retires the now-unused "2017-11-20 for 3.30.09.00" shared string.
$ws.Range("F1").Value = "2018-01-25 for 3.30.10.01"

# --- Row 48: 3.30.10.01 / new / read / fix logic error in settlement age calc ---
$ws.Range("A48").Value = 43125
$ws.Range("B48").Value = "3.30.10.01"
$ws.Range("C48").Value = "new"
$ws.Range("D48").Value = "read"
$ws.Range("F48").Value = "fix logic error in calculation of settlement age in 2 season, 2 settlement setup"
$ws.Range("G48").Value = "No"

# Move selection to reflect the newly active editing location.
$ws.Activate() | Out-Null
$ws.Range("F50").Select() | Out-Null
